# Apply cryptos list update (values refreshed by scheduled scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.428.57"
$ws.Range("E2").Value = "  -3.97%  "
$ws.Range("D3").Value = "2.458.96"
$ws.Range("E3").Value = "  -6.83%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.49%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.91%  "
$ws.Range("D9").Value = "2.457.23"
$ws.Range("E9").Value = "  -6.86%  "
$ws.Range("E10").Value = "  -9.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.42%  "
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.87%  "
$ws.Range("D15").Value = "2.905.10"
$ws.Range("E15").Value = "  -6.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -10.14%  "
$ws.Range("D17").Value = "61.360.18"
$ws.Range("E17").Value = "  -3.93%  "
$ws.Range("D18").Value = "2.462.35"
$ws.Range("E18").Value = "  -6.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "318.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.63%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.96%  "
$ws.Range("D26").Value = "0.0₃0988"
$ws.Range("E26").Value = "  -12.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "560.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.62%  "
$ws.Range("D28").Value = "2.602.74"
$ws.Range("E28").Value = "  -6.11%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.148"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.96%  "
$ws.Range("E34").Value = "  -6.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.17%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.379"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "142.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.37%  "
$ws.Range("E42").Value = "  -8.83%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0536"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.591"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0940"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.30%  "
